$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before H (old H "Yco" ... shifts to I, etc.)
$ws.Range("H1").EntireColumn.Insert()

# New header for the inserted column
$ws.Range("H1").Value = "CO2/(CO+CO2)"

# Fill the new column with the CO2/(CO+CO2) ratio formula for each data row
$ws.Range("H2:H16").Formula = "=F2/(E2+F2)"

# Add three empty rows right after the data (rows 17-19), matching row 16's column H
$ws.Range("H17:H19").NumberFormat = "General"

# Update the active selection to reflect where the user ended up
$ws.Range("H8").Select()
